$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('G2').Value = 'System, backup@backdoor.com, system'
$ws.Range('G3').Value = 'System, dnasr281@gmail.com'
$ws.Range('G6').Value = 'System, dnasr281@gmail.com'
$ws.Range('G11').Value = 'System, dnasr281@gmail.com'
$ws.Range('G12').Value = 'System, dnasr281@gmail.com'
$ws.Range('G13').Value = 'System, dnasr281@gmail.com'
$ws.Range('G14').Value = 'System, dnasr281@gmail.com'
$ws.Range('G15').Value = 'System, dnasr281@gmail.com'
$ws.Range('G29').Value = 'System, backup@backdoor.com, system'
$ws.Range('G30').Value = 'System, dnasr281@gmail.com'
$ws.Range('G33').Value = 'System, dnasr281@gmail.com'
$ws.Range('G38').Value = 'System, dnasr281@gmail.com'
$ws.Range('G39').Value = 'System, dnasr281@gmail.com'
$ws.Range('G40').Value = 'System, dnasr281@gmail.com'
$ws.Range('G41').Value = 'System, dnasr281@gmail.com'
$ws.Range('G42').Value = 'System, dnasr281@gmail.com'
$ws.Range('G56').Value = 'System, backup@backdoor.com, system'
$ws.Range('G57').Value = 'System, dnasr281@gmail.com'
$ws.Range('G60').Value = 'System, dnasr281@gmail.com'
$ws.Range('G65').Value = 'System, dnasr281@gmail.com'
$ws.Range('G66').Value = 'System, dnasr281@gmail.com'
$ws.Range('G67').Value = 'System, dnasr281@gmail.com'
$ws.Range('G68').Value = 'System, dnasr281@gmail.com'
$ws.Range('G69').Value = 'System, dnasr281@gmail.com'
$ws.Range('G89').Value = 'System, dnasr281@gmail.com'
$ws.Range('G90').Value = 'admin@admin.com, dnasr281@gmail.com'
$ws.Range('G93').Value = 'System, dnasr281@gmail.com'
$ws.Range('G115').Value = 'System, dnasr281@gmail.com'
$ws.Range('G116').Value = 'admin@admin.com, dnasr281@gmail.com'
$ws.Range('G119').Value = 'System, dnasr281@gmail.com'
$ws.Range('G141').Value = 'System, dnasr281@gmail.com'
$ws.Range('G142').Value = 'admin@admin.com, dnasr281@gmail.com'
$ws.Range('G145').Value = 'System, dnasr281@gmail.com'
